# "Adds all IG authors as contact"
#
# The Metadata sheet lists one "Contact" row per IG author (rows 10 and 11
# already hold a "Contact" / "No display for ContactDetail" pair each).
# Two more authors were added as contacts, so we insert two more identical
# "Contact" rows directly below the existing ones (at rows 12-13), which
# pushes every row below it (Jurisdiction, Description, Purpose, ...) down
# by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert two blank rows right after the last existing "Contact" row (row 11).
$ws.Rows.Item(12).Resize(2).Insert()

# Copy the formatting (borders/fill/alignment/style) from the row above so the
# newly inserted rows look identical to the other data rows instead of
# picking up a brand new, unformatted style.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new "Contact" rows with the same text used for the existing
# contact entries.
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"
